$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Cód barras" column (T) used for the product alias shown only on PDFs/quotations
$ws.Range("T1").Value = "Cód barras"
$ws.Range("T2").Value = 10000001
$ws.Range("T3").Value = 10000002

# Underline an (otherwise empty) far cell, matching the formatting probe left behind in the edit session
$ws.Range("S8").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle

# Restore cursor / viewport position left after the edits
$ws.Range("S8").Select()
$excel.ActiveWindow.ScrollColumn = 12
